$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 174; existing rows 174.. shift down to 175..
$ws.Rows.Item(174).Insert()

# Populate the newly inserted row 174 with this week's new record
$ws.Cells.Item(174, 1).Value = 11
$ws.Cells.Item(174, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(174, 3).Value = "Bíobío"
$ws.Cells.Item(174, 4).Value = 45029
$ws.Cells.Item(174, 5).Value = 8
$ws.Cells.Item(174, 6).Value = "Fruta"
$ws.Cells.Item(174, 7).Value = 100108
$ws.Cells.Item(174, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(174, 9).Value = 100108005
$ws.Cells.Item(174, 10).Value = "Piña"
$ws.Cells.Item(174, 11).Value = "Caramelo"
$ws.Cells.Item(174, 12).Value = "Segunda"
$ws.Cells.Item(174, 13).Value = 100
$ws.Cells.Item(174, 14).Value = 19000
$ws.Cells.Item(174, 15).Value = 20000
$ws.Cells.Item(174, 16).Value = 19500
$ws.Cells.Item(174, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item(174, 18).Value = "Ecuador"
$ws.Cells.Item(174, 19).Value = 1393
$ws.Cells.Item(174, 20).Value = 14

# Keep the date column formatted the same way as the rest of column D
$ws.Cells.Item(174, 4).NumberFormat = $ws.Cells.Item(175, 4).NumberFormat
